$d = $word.ActiveDocument

# The site footer (the "Ver no Jupiter..." line, the copyright line, and the
# blank paragraph that separates them from the body text) was dropped from
# the rebuilt page, so remove those trailing paragraphs from the document.

function Find-ParagraphIndex($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

$jupiterIdx = Find-ParagraphIndex "Ver no Jupiter Salvar em pdf Salvar em docx"
$copyrightIdx = Find-ParagraphIndex "Contact: luizeleno@usp.br"

if ($jupiterIdx -gt 0 -and $copyrightIdx -gt 0) {
    # Also drop the blank paragraph immediately preceding the "Ver no
    # Jupiter..." paragraph, so the Requisitos text is followed directly by
    # the (still-present) trailing blank paragraph / page break.
    $blankIdx = $jupiterIdx - 1
    $startPara = $d.Paragraphs.Item($blankIdx)
    $endPara = $d.Paragraphs.Item($copyrightIdx)

    $delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $delRange.Delete() | Out-Null
}

Write-Host "Paragraphs remaining:" $d.Paragraphs.Count
